# Updated cryptos list with refreshed price/volume data (and a couple of
# coin rows that changed ranking position) as scraped from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.247.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.58%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.376.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.76%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'540.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.70%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'140.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.05%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.14%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.61%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.376.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.71%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.46%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +0.12%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'5.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.86%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.03%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'25.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.82%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.802.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.91%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.31%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'60.213.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.46%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.372.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.75%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'10.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.56%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "'BitcoinCash"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'315.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.65%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = "'Polkadot"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'4.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.58%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.64%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.06%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.23%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'62.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.22%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.40%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.498.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.71%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.0₃0929"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.32%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.41%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'Fetch.AI"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.83%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'Bittensor"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'512.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.93%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.22%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -4.51%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -3.55%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.08%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.20%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'4.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.95%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'5.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.54%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.41%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'17.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.45%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D42").Value = "'1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.23%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'137.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.68%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'40.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.18%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.27%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'139.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -6.15%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.27%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'20.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.06%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.46%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.574"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.36%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0922"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.16%  "
$ws.Range("E51").Style = "Normal"
